$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update row 2 values (B2:E2)
$ws.Range("B2").Value = 11.310159249463084
$ws.Range("C2").Value = -0.60162695947590805
$ws.Range("D2").Value = 0.041881807647570213
$ws.Range("E2").Value = -1.2772184283606407

# Update row 3 values (B3:E3)
$ws.Range("B3").Value = 4.0596357921674269
$ws.Range("C3").Value = 4.9789892904264548
$ws.Range("D3").Value = 1.231805422257537
$ws.Range("E3").Value = -3.9191652716807255

# Update selection range to match new sqref B1:E3
$ws.Range("B1:E3").Select()
